$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 418, shifting existing rows 418-468 down to 419-469.
$ws.Rows.Item(418).Insert()

# Populate the newly inserted row 418 with the new data point.
$ws.Cells.Item(418, 1).Value = 5
$ws.Cells.Item(418, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(418, 3).Value = "Maule"
$ws.Cells.Item(418, 4).Value = 44578
$ws.Cells.Item(418, 5).Value = 7
$ws.Cells.Item(418, 6).Value = "Fruta"
$ws.Cells.Item(418, 7).Value = 100104
$ws.Cells.Item(418, 8).Value = "Frutos de pepita"
$ws.Cells.Item(418, 9).Value = 100104005
$ws.Cells.Item(418, 10).Value = "Pera"
$ws.Cells.Item(418, 11).Value = "Carmen"
$ws.Cells.Item(418, 12).Value = "Segunda"
$ws.Cells.Item(418, 13).Value = 300
$ws.Cells.Item(418, 14).Value = 9000
$ws.Cells.Item(418, 15).Value = 9000
$ws.Cells.Item(418, 16).Value = 9000
$ws.Cells.Item(418, 17).Value = "`$/bandeja 18 kilos granel"
$ws.Cells.Item(418, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(418, 19).Value = 500
$ws.Cells.Item(418, 20).Value = 18
